$wb = $excel.ActiveWorkbook

$wsManu = $wb.Worksheets.Item("Manufacturer")
$wsManu.Range("D7").Value = "quantity"

$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRetailer = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$wsRetailer.Name = "Retailer"

# Copy formats from Manufacturer sheet rows that already have the right styles
$wsManu.Range("A1:B1").Copy()
$wsRetailer.Range("A1:B1").PasteSpecial(-4122)

$wsManu.Range("A2:B2").Copy()
$wsRetailer.Range("A2:B2").PasteSpecial(-4122)

$wsManu.Range("A7:D7").Copy()
$wsRetailer.Range("A4:D4").PasteSpecial(-4122)

$wsManu.Range("A8:D8").Copy()
$wsRetailer.Range("A5:D5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Now set values
$wsRetailer.Range("A1").Value = "Sl No"
$wsRetailer.Range("B1").Value = "TestCaseName"
$wsRetailer.Range("A2").Value = 1
$wsRetailer.Range("B2").Value = "retailerLoginTest"

$wsRetailer.Range("A4").Value = "Sl No"
$wsRetailer.Range("B4").Value = "TestCaseName"
$wsRetailer.Range("C4").Value = "productName"
$wsRetailer.Range("D4").Value = "quantity"
$wsRetailer.Range("A5").Value = 2
$wsRetailer.Range("B5").Value = "createAOrder"
$wsRetailer.Range("C5").Value = "pagent"
$wsRetailer.Range("D5").Value = "'143"

$wsManu.Range("B5").Value = "manufacturerAddProductTest"
$wsManu.Range("B2").Value = "manufacturerLoginTest"

# Set selections: Manufacturer -> B2, Retailer -> A4:D5 (activeCell A4), Retailer active last
$null = $wsManu.Range("B2").Select()
$null = $wsRetailer.Range("A4:D5").Select()

Write-Output "done"
